$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.89420086145401
$ws.Range("B1").Value = 2.951676607131958
$ws.Range("C1").Value = 7.009454727172852
$ws.Range("D1").Value = 2.019955635070801
$ws.Range("E1").Value = 1.406994462013245
